# Apply updated cryptocurrency price/volume data to the "cryptos" sheet.
# Price (D) and Volume(1h) (E) are stored as literal text (e.g. "30.411.46",
# "1.370", "  -1.15%  ") -- NOT numbers/percentages -- so we force the
# NumberFormat to Text ("@") before writing, otherwise Excel's COM layer
# auto-converts numeric-looking strings (losing trailing zeros, e.g.
# "1.370" -> 1.37, "53.00" -> 53) and percent-looking strings into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '30.411.46'
$ws.Range('E2').Value = '  -1.15%  '
$ws.Range('D3').Value = '1.917.35'
$ws.Range('E3').Value = '  +1.79%  '
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Value = '240.78'
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('D7').Value = '0.4683'
$ws.Range('E7').Value = '  -2.03%  '
$ws.Range('D8').Value = '0.2846'
$ws.Range('E8').Value = '  +0.27%  '
$ws.Range('D9').Value = '0.06831'
$ws.Range('E9').Value = '  +5.27%  '
$ws.Range('D10').Value = '107.03'
$ws.Range('E10').Value = '  +13.55%  '
$ws.Range('D11').Value = '18.04'
$ws.Range('E11').Value = '  -4.28%  '
$ws.Range('D12').Value = '1.907.76'
$ws.Range('E12').Value = '  +1.22%  '
$ws.Range('D13').Value = '0.07626'
$ws.Range('E13').Value = '  +0.89%  '
$ws.Range('D14').Value = '5.179'
$ws.Range('E14').Value = '  +1.62%  '
$ws.Range('D15').Value = '0.6547'
$ws.Range('E15').Value = '  +0.67%  '
$ws.Range('D16').Value = '285.33'
$ws.Range('E16').Value = '  -3.58%  '
$ws.Range('D17').Value = '30.373.72'
$ws.Range('E17').Value = '  -1.24%  '
$ws.Range('D18').Value = '0.9999'
$ws.Range('E18').Value = '  -0.23%  '
$ws.Range('D19').Value = '0.000007592'
$ws.Range('E19').Value = '  +1.51%  '
$ws.Range('D20').Value = '12.95'
$ws.Range('E20').Value = '  -1.38%  '
$ws.Range('D21').Value = '2.157.17'
$ws.Range('E21').Value = '  +0.81%  '
$ws.Range('D22').Value = '0.9987'
$ws.Range('E22').Value = '  -0.39%  '
$ws.Range('D23').Value = '5.209'
$ws.Range('E23').Value = '  +1.37%  '
$ws.Range('D24').Value = '6.195'
$ws.Range('E24').Value = '  +1.02%  '
$ws.Range('D25').Value = '168.04'
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('D26').Value = '9.241'
$ws.Range('E26').Value = '  -0.42%  '
$ws.Range('D27').Value = '21.39'
$ws.Range('E27').Value = '  +9.52%  '
$ws.Range('D28').Value = '2.041'
$ws.Range('E28').Value = '  +4.80%  '
$ws.Range('D29').Value = '0.1071'
$ws.Range('E29').Value = '  +1.02%  '
$ws.Range('D30').Value = '1.370'
$ws.Range('E30').Value = '  +1.00%  '
$ws.Range('D31').Value = '4.131'
$ws.Range('E31').Value = '  -1.15%  '
$ws.Range('D32').Value = '3.942'
$ws.Range('E32').Value = '  -0.27%  '
$ws.Range('D33').Value = '0.05033'
$ws.Range('E33').Value = '  +0.24%  '
$ws.Range('D34').Value = '0.7358'
$ws.Range('E34').Value = '  +2.17%  '
$ws.Range('D35').Value = '1.144'
$ws.Range('E35').Value = '  -2.09%  '
$ws.Range('D36').Value = '0.9992'
$ws.Range('D37').Value = '2.729'
$ws.Range('E37').Value = '  +0.48%  '
$ws.Range('D38').Value = '0.02019'
$ws.Range('E38').Value = '  +3.93%  '
$ws.Range('D39').Value = '2.680'
$ws.Range('E39').Value = '  -1.19%  '
$ws.Range('D40').Value = '2.047'
$ws.Range('E40').Value = '  -0.41%  '
$ws.Range('D41').Value = '108.76'
$ws.Range('E41').Value = '  +1.64%  '
$ws.Range('D42').Value = '0.8716'
$ws.Range('E42').Value = '  -2.76%  '
$ws.Range('D43').Value = '5.812'
$ws.Range('E43').Value = '  +4.28%  '

# BitcoinSV and PaxDollar swapped rank positions (rows 44/45 fully swapped)
$ws.Range('B44').Value = 'BitcoinSV'
$ws.Range('C44').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D44').Value = '53.00'
$ws.Range('E44').Value = '  +26.27%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').Value = '0.9997'
$ws.Range('E45').Value = '  -0.21%  '
$ws.Range('D46').Value = '0.4194'
$ws.Range('E46').Value = '  +0.07%  '
$ws.Range('D47').Value = '67.30'
$ws.Range('E47').Value = '  +2.02%  '
$ws.Range('D48').Value = '7.146'
$ws.Range('E48').Value = '  -2.43%  '
$ws.Range('D49').Value = '9.192'
$ws.Range('E49').Value = '  +4.39%  '
$ws.Range('D50').Value = '0.1206'
$ws.Range('E50').Value = '  -1.23%  '
$ws.Range('D51').Value = '34.57'
$ws.Range('E51').Value = '  -0.02%  '
